$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column B (Probabilidad Simulada) values for rows 3-95 per Monte Carlo
# simulation re-run (commit: "Se agrega prueba de modelos")
$ws.Cells.Item(3, 2).Value = 0.06610000000000001
$ws.Cells.Item(4, 2).Value = 0.0475
$ws.Cells.Item(5, 2).Value = 0.0766
$ws.Cells.Item(6, 2).Value = 0.0882
$ws.Cells.Item(7, 2).Value = 0.1074
$ws.Cells.Item(8, 2).Value = 0.1002
$ws.Cells.Item(9, 2).Value = 0.0979
$ws.Cells.Item(10, 2).Value = 0.1036
$ws.Cells.Item(11, 2).Value = 0.0941
$ws.Cells.Item(12, 2).Value = 0.0762
$ws.Cells.Item(13, 2).Value = 0.0634
$ws.Cells.Item(14, 2).Value = 0.0788
$ws.Cells.Item(15, 2).Value = 0.0843
$ws.Cells.Item(16, 2).Value = 0.0682
$ws.Cells.Item(17, 2).Value = 0.0467
$ws.Cells.Item(18, 2).Value = 0.1112
$ws.Cells.Item(19, 2).Value = 0.121
$ws.Cells.Item(20, 2).Value = 0.1101
$ws.Cells.Item(21, 2).Value = 0.0901
$ws.Cells.Item(22, 2).Value = 0.0956
$ws.Cells.Item(23, 2).Value = 0.0843
$ws.Cells.Item(24, 2).Value = 0.0672
$ws.Cells.Item(25, 2).Value = 0.07190000000000001
$ws.Cells.Item(26, 2).Value = 0.0494
$ws.Cells.Item(27, 2).Value = 0.0858
$ws.Cells.Item(28, 2).Value = 0.0577
$ws.Cells.Item(29, 2).Value = 0.09950000000000001
$ws.Cells.Item(30, 2).Value = 0.09660000000000001
$ws.Cells.Item(31, 2).Value = 0.1074
$ws.Cells.Item(32, 2).Value = 0.0997
$ws.Cells.Item(33, 2).Value = 0.0717
$ws.Cells.Item(34, 2).Value = 0.0948
$ws.Cells.Item(35, 2).Value = 0.0838
$ws.Cells.Item(36, 2).Value = 0.066
$ws.Cells.Item(37, 2).Value = 0.07920000000000001
$ws.Cells.Item(38, 2).Value = 0.0578
$ws.Cells.Item(39, 2).Value = 0.0805
$ws.Cells.Item(40, 2).Value = 0.06270000000000001
$ws.Cells.Item(41, 2).Value = 0.064
$ws.Cells.Item(42, 2).Value = 0.1165
$ws.Cells.Item(43, 2).Value = 0.1007
$ws.Cells.Item(44, 2).Value = 0.1259
$ws.Cells.Item(45, 2).Value = 0.0755
$ws.Cells.Item(46, 2).Value = 0.0663
$ws.Cells.Item(47, 2).Value = 0.08550000000000001
$ws.Cells.Item(48, 2).Value = 0.09
$ws.Cells.Item(49, 2).Value = 0.0786
$ws.Cells.Item(50, 2).Value = 0.0538
$ws.Cells.Item(51, 2).Value = 0.0693
$ws.Cells.Item(52, 2).Value = 0.0555
$ws.Cells.Item(53, 2).Value = 0.1024
$ws.Cells.Item(54, 2).Value = 0.0914
$ws.Cells.Item(55, 2).Value = 0.0838
$ws.Cells.Item(56, 2).Value = 0.092
$ws.Cells.Item(57, 2).Value = 0.0897
$ws.Cells.Item(58, 2).Value = 0.07729999999999999
$ws.Cells.Item(59, 2).Value = 0.09619999999999999
$ws.Cells.Item(60, 2).Value = 0.067
$ws.Cells.Item(61, 2).Value = 0.0982
$ws.Cells.Item(62, 2).Value = 0.0772
$ws.Cells.Item(63, 2).Value = 0.0535
$ws.Cells.Item(64, 2).Value = 0.0539
$ws.Cells.Item(65, 2).Value = 0.0882
$ws.Cells.Item(66, 2).Value = 0.0921
$ws.Cells.Item(67, 2).Value = 0.1114
$ws.Cells.Item(68, 2).Value = 0.1083
$ws.Cells.Item(69, 2).Value = 0.08169999999999999
$ws.Cells.Item(70, 2).Value = 0.1072
$ws.Cells.Item(71, 2).Value = 0.0814
$ws.Cells.Item(72, 2).Value = 0.0953
$ws.Cells.Item(73, 2).Value = 0.068
$ws.Cells.Item(74, 2).Value = 0.059
$ws.Cells.Item(75, 2).Value = 0.08069999999999999
$ws.Cells.Item(76, 2).Value = 0.0629
$ws.Cells.Item(77, 2).Value = 0.1002
$ws.Cells.Item(78, 2).Value = 0.08119999999999999
$ws.Cells.Item(79, 2).Value = 0.1063
$ws.Cells.Item(80, 2).Value = 0.1133
$ws.Cells.Item(81, 2).Value = 0.0756
$ws.Cells.Item(82, 2).Value = 0.0737
$ws.Cells.Item(83, 2).Value = 0.0825
$ws.Cells.Item(84, 2).Value = 0.0813
$ws.Cells.Item(85, 2).Value = 0.0589
$ws.Cells.Item(86, 2).Value = 0.0834
$ws.Cells.Item(87, 2).Value = 0.09130000000000001
$ws.Cells.Item(88, 2).Value = 0.1198
$ws.Cells.Item(89, 2).Value = 0.0958
$ws.Cells.Item(90, 2).Value = 0.1375
$ws.Cells.Item(91, 2).Value = 0.1395
$ws.Cells.Item(92, 2).Value = 0.1082
$ws.Cells.Item(93, 2).Value = 0.09379999999999999
$ws.Cells.Item(94, 2).Value = 0.0927
$ws.Cells.Item(95, 2).Value = 0.1214
